$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "301.80"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.97%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "37.47"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "6.28%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-2.76%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07855"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.02%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.231"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-7.52%"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.15%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "4.019"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "2.05%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9081"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.61%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1889"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "4.94%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09424"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-5.63%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08507"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.02%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03522"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "6.32%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09964"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.73%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001481"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.16%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005716"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.45%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.468"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.06%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.077"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-2.93%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2.94%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.55%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.768"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "10.69%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2202"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-7.62%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04641"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.56%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001228"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.04%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004454"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.09%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001299"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.00%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0004746"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "28.45%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01771"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-1.21%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04749"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.04%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007799"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.79%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1391"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-1.45%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007658"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "7.48%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002228"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "5.89%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009823"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "3.18%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006040"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-1.18%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000749"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.05%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "217.35%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002688"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "34.57%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002098"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.05%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001998"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.05%"
